# Update column G ("K") values on Sheet1 with regenerated strikeout counts.
# (commit: "regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    2  = 6
    3  = 7
    4  = 6
    5  = 8
    6  = 0
    7  = 6
    8  = 3
    9  = 6
    10 = 0
    11 = 7
    12 = 4
    13 = 4
    14 = 7
    15 = 6
    16 = 5
    17 = 9
    18 = 7
    19 = 4
    20 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
